# Add Q4-2022 fund-holding data for 300795-米奥会展
#
# 1. Insert a new worksheet named "2022-Q4" right after "总计" (i.e. before
#    the current first quarter sheet, which shifts "2022-Q3", "2022-Q2",
#    "2021-Q1" one position to the right).
# 2. Populate "2022-Q4" with the same column layout used by the other
#    quarter sheets (基金代码/基金名称/基金规模/股票总仓位/仓位占比/
#    持有市值(亿元)/仓位排名) and copy header/index-column formatting from
#    the existing "2022-Q3" sheet so the new sheet matches house style.
# 3. Update the "总计" summary sheet: insert a new row for 2022-Q4 right
#    below the header (count = number of funds, value = sum of 持有市值)
#    and keep the existing quarters below it, renumbering the index column.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2022-Q3")

# ---------------------------------------------------------------------
# Step 1: create the new sheet before "2022-Q3" and name it "2022-Q4"
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($q3, $null)
$q4.Name = "2022-Q4"

# ---------------------------------------------------------------------
# Step 2: copy formatting (bold/border/center) from "2022-Q3" header row
# and index column onto the same cells of the new "2022-Q4" sheet, then
# fill in the values.
# ---------------------------------------------------------------------
$q3.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

$q3.Range("A2").Copy()
$q4.Range("A2:A8").PasteSpecial(-4122)

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4Data = @(
    @(0, "003292", "嘉实优势成长灵活配置混合", "8.43", "87.02", "4.96", "0.4181", 4),
    @(1, "070022", "嘉实领先成长混合", "4.88", "90.23", "5.03", "0.2455", 5),
    @(2, "009381", "汇安核心资产混合A", "3.84", "92.63", "4.91", "0.1885", 5),
    @(3, "160722", "嘉实惠泽灵活配置混合（LOF）", "0.87", "92.18", "8.23", "0.0716", 1),
    @(4, "007775", "汇安量化先锋混合A", "0.22", "86.70", "4.92", "0.0108", 2),
    @(5, "007776", "汇安量化先锋混合C", "0.12", "86.70", "4.92", "0.0059", 2),
    @(6, "009382", "汇安核心资产混合C", "0.07", "92.63", "4.91", "0.0034", 5)
)

$r = 2
foreach ($row in $q4Data) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = "'" + $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = "'" + $row[3]
    $q4.Cells.Item($r, 5).Value = "'" + $row[4]
    $q4.Cells.Item($r, 6).Value = "'" + $row[5]
    $q4.Cells.Item($r, 7).Value = "'" + $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 3: update the "总计" (summary) sheet with a new row for 2022-Q4.
# ---------------------------------------------------------------------
$totalRows = @(
    @("2022-Q4", 7, 0.9399999999999999),
    @("2022-Q3", 4, 0.42),
    @("2022-Q2", 3, 0.35),
    @("2021-Q1", 2, 0)
)

$total.Range("A2").Copy()
$total.Range("A2:A5").PasteSpecial(-4122)

$r = 2
$idx = 0
foreach ($row in $totalRows) {
    $total.Cells.Item($r, 1).Value = $idx
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
    $idx = $idx + 1
}

# Restore original active sheet / selection state.
$total.Activate()
$total.Range("A1").Select() | Out-Null
